$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values A8:A12 from 100 to 200
$ws.Range("A8").Value = 200
$ws.Range("A9").Value = 200
$ws.Range("A10").Value = 200
$ws.Range("A11").Value = 200
$ws.Range("A12").Value = 200

# Update the active cell/selection to A13
$ws.Range("A13").Select()
